$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.558.76'
$ws.Range("E2").Value = '  -4.64%  '
$ws.Range("D3").Value = '2.936.94'
$ws.Range("E3").Value = '  -2.49%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '548.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.52%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.23'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.89%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.510'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.22%  '
$ws.Range("D9").Value = '2.929.73'
$ws.Range("E9").Value = '  -2.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.127'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.45%  '
$ws.Range("E11").Value = '  -6.07%  '
$ws.Range("E12").Value = '  +0.75%  '
$ws.Range("E13").Value = '  -0.24%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.87'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.52%  '
$ws.Range("E15").Value = '  +0.78%  '
$ws.Range("D16").Value = '3.422.37'
$ws.Range("E16").Value = '  -2.45%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.89'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +6.26%  '
$ws.Range("D18").Value = '2.932.92'
$ws.Range("E18").Value = '  -2.35%  '
$ws.Range("D19").Value = '57.604.67'
$ws.Range("E19").Value = '  -4.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '417.34'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.25'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.60%  '
$ws.Range("E22").Value = '  +2.38%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.99'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.43%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.09'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '79.74'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.33%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("E28").Value = '  -3.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.45'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.38%  '
$ws.Range("E30").Value = '  +0.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '25.16'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.88%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.01'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.89%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0968'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.96%  '
$ws.Range("E34").Value = '  +0.84%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.937'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.29%  '
$ws.Range("E36").Value = '  +0.56%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '48.28'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.99%  '
$ws.Range("D38").Value = '0.0₃0683'
$ws.Range("E38").Value = '  +1.12%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.65'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.30%  '
$ws.Range("E40").Value = '  +3.17%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.108'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.29%  '
$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '376.20'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("E43").Value = '  -3.41%  '
$ws.Range("D44").Value = '2.695.71'
$ws.Range("E44").Value = '  +0.50%  '
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("E46").Value = '  +0.96%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '122.36'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.90%  '
$ws.Range("E48").Value = '  +1.16%  '
$ws.Range("E49").Value = '  -2.20%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.08'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.00%  '
$ws.Range("E51").Value = '  -0.20%  '
